$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.924.33"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.34%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.908.98"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -3.13%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'591.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -0.22%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'145.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -3.95%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.13%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.61%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.907.36"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -2.78%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'6.72"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.20%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -2.71%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  -3.28%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  -1.30%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'33.56"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -4.58%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  +0.26%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.392.40"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -3.15%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'60.853.47"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.69%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'6.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -3.88%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'2.908.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -3.33%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'429.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -3.29%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'13.53"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -4.02%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  -1.03%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'7.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -4.80%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'81.32"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.24%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'10.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.94%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'2.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -1.28%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'11.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.91%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +0.07%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.32"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +2.97%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  -0.28%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'2.62"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -2.88%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  -4.87%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  -2.68%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  -1.91%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.0₃0854"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.31%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -2.26%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'5.61"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -3.74%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'3.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -1.68%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'49.58"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -1.46%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.124"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -2.60%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'2.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -3.61%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'8.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.52%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.291"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -3.58%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'40.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -9.28%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'375.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.43%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.0348"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -2.46%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'2.706.92"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.54%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'129.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -2.60%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'24.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -7.91%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.106"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.43%  "
$ws.Range("E51").ClearFormats()
